$wb = $excel.ActiveWorkbook

# Sheet 1: Estadistica_general - R_deming row (row 3), column B
$ws1 = $wb.Worksheets.Item("Estadistica_general")
$ws1.Range("B3").Value = 0.2954342378901791

# Sheet 2: Estadistica_mensual - R_deming row (row 3), columns E:I
$ws2 = $wb.Worksheets.Item("Estadistica_mensual")
$ws2.Range("E3").Value = 1.655896483460473
$ws2.Range("F3").Value = 0.8076142579283497
$ws2.Range("G3").Value = -0.540900480405952
$ws2.Range("H3").Value = 0.2358961189181862
$ws2.Range("I3").Value = -0.3263886678766575

# Sheet 3: Estadistica_anual - R_deming row (row 3), columns B:D
$ws3 = $wb.Worksheets.Item("Estadistica_anual")
$ws3.Range("B3").Value = 0.2862379893085237
$ws3.Range("C3").Value = 0.2457007524798652
$ws3.Range("D3").Value = 0.4717234119244378
